$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 14/15 previously had an empty "time" cell (merged L:M); fill them in
# with "45 phút " (trailing space preserved from the source data).
$ws.Range("L14").Value = "45 phút "
$ws.Range("L15").Value = "45 phút "

# Rows 17/18 previously said "30 phút"; update to "45 phút".
$ws.Range("L17").Value = "45 phút"
$ws.Range("L18").Value = "45 phút"

# Move the active selection to L19 (was O19).
$ws.Range("L19").Select()
